$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 13: Horai Gakuen DX1 Revised Replay Special (SoftBank Creative)
$ws.Range("B13").Value2 = '蓬莱学園DX1 改訂版リプレイ特集'
$ws.Range("C13").Value2 = 'Horai Gakuen DX1 Revised Replay Special'
# Row 14: Horai Gakuen DX2 Special Feature / 1995 School Festival
$ws.Range("B14").Value2 = '蓬莱学園DX2 特集・1995年学園祭'
$ws.Range("C14").Value2 = 'Horai Gakuen DX2 Special Feature / 1995 School Festival'
# Row 13: Horai Gakuen DX1 Revised Replay Special (SoftBank Creative)
$ws.Range("D13").Value2 = 'SoftBank Creative'
# Row 15: Horai Gakuen DX3 Class / Graduation
$ws.Range("B15").Value2 = '蓬莱学園DX3 蓬莱学園の授業・卒業'
$ws.Range("C15").Value2 = 'Horai Gakuen DX3 Horai Gakuen Class / Graduation'
# Row 13: Horai Gakuen DX1 Revised Replay Special (SoftBank Creative)
$ws.Range("E13").Value2 = 'revised-replay.jpg'
$ws.Range("F13").Value2 = 'replay'
$ws.Range("A13").Value2 = 1995
# Row 14: Horai Gakuen DX2 Special Feature / 1995 School Festival
$ws.Range("E14").Value2 = 'revised-school-feature.jpg'
$ws.Range("A14").Value2 = 1995
$ws.Range("D14").Value2 = 'SoftBank Creative'
$ws.Range("F14").Value2 = 'supplement'
# Row 15: Horai Gakuen DX3 Class / Graduation
$ws.Range("E15").Value2 = 'revised-graduation.jpg'
$ws.Range("A15").Value2 = 1996
$ws.Range("D15").Value2 = 'SoftBank Creative'
$ws.Range("F15").Value2 = 'supplement'
# Row 12: Horai Gakuen World Tour (Shinkigensha)
$ws.Range("B12").Value2 = '蓬莱学園ワールドツアー'
$ws.Range("C12").Value2 = 'Horai Gakuen World Tour'
$ws.Range("D12").Value2 = 'Shinkigensha'
$ws.Range("E12").Value2 = 'hourai_gakuen_world_tour.jpg'
$ws.Range("A12").Value2 = 1994
$ws.Range("F12").Value2 = 'supplement'
# Row 8: Horai Life Encyclopedia (Yuentai)
$ws.Range("B8").Value2 = '蓬莱生活事典　蓬莱学園の冒険！！ー改訂版ーサプリメント'
$ws.Range("C8").Value2 = 'Horai Life Encyclopedia: Adventure of Horai Gakuen! !! Revised version'
$ws.Range("E8").Value2 = 'horai_school_living_encyclopdedia.jpg'
$ws.Range("D8").Value2 = 'Yuentai'
$ws.Range("F8").Value2 = 'supplement'
# Row 9: Horai Gakuen Biographical Dictionary (Yuentai)
$ws.Range("B9").Value2 = '蓬莱学園人名事典'
$ws.Range("C9").Value2 = 'Horai Gakuen Biographical Dictionary'
$ws.Range("E9").Value2 = 'horai_gakuen_biographical_dictionary.jpg'
$ws.Range("A9").Value2 = 1995
$ws.Range("D9").Value2 = 'Yuentai'
$ws.Range("F9").Value2 = 'supplement'
# Row 11: Everything about Horai Gakuen (BNN)
$ws.Range("B11").Value2 = 'なんでもかんでも蓬莱学園'
# Row 10: Taking the Exam Horai Gakuen (BNN)
$ws.Range("D10").Value2 = 'BNN'
# Row 11: Everything about Horai Gakuen (BNN)
$ws.Range("C11").Value2 = 'Everything about Horai Gakuen'
$ws.Range("E11").Value2 = 'everything_about_horai_gakuen.jpg'
$ws.Range("A11").Value2 = 1994
$ws.Range("D11").Value2 = 'BNN'
$ws.Range("F11").Value2 = 'supplement'
# Row 10: Taking the Exam Horai Gakuen (BNN)
$ws.Range("B10").Value2 = '試験に出る蓬莱学園!'
$ws.Range("C10").Value2 = 'Taking the Exam Horai Gakuen'
$ws.Range("E10").Value2 = 'taking_the_exam.jpg'
$ws.Range("A10").Value2 = 1991
$ws.Range("F10").Value2 = 'supplement'
# New 'product_code' column (G)
$ws.Range("G1").Value2 = 'product_code'
$ws.Range("G2").Value2 = 'Y-9101'
$ws.Range("G5").Value2 = 'Y-9104'
$ws.Range("A5").Value2 = 1993

# Row 5 gains a publication year
$ws.Range("A5").Value2 = 1993

# Columns B and C widened to fit the new, longer titles
$ws.Columns.Item(2).ColumnWidth = 58.0
$ws.Columns.Item(3).ColumnWidth = 62.5

# Selection moved to A6
$ws.Range("A6").Select() | Out-Null
